# Apply per-cell value updates for the team-specific time matrix (Colgate_A)
# as described in the commit: "added team specific time data, have not yet
# implemented its logic for simulation". Only numeric <v> contents change;
# no styles, formulas, or structure are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2193308550185874
$ws.Range("C2").Value = 0.5018587360594795
$ws.Range("J2").Value = 0.01858736059479554
$ws.Range("P2").Value = 0.1524163568773234
$ws.Range("S2").Value = 0.1078066914498141

# Row 3
$ws.Range("B3").Value = 0.007352941176470588
$ws.Range("C3").Value = 0.01470588235294118
$ws.Range("J3").Value = 0.05882352941176471
$ws.Range("P3").Value = 0.6985294117647058
$ws.Range("S3").Value = 0.2205882352941176

# Row 4
$ws.Range("P4").Value = 0.7555555555555555
$ws.Range("S4").Value = 0.2444444444444444

# Row 5
$ws.Range("J5").Value = 0.4
$ws.Range("P5").Value = 0.4
$ws.Range("S5").Value = 0.2

# Row 6
$ws.Range("B6").Value = 0.06572769953051644
$ws.Range("D6").Value = 0.0187793427230047
$ws.Range("E6").Value = 0.004694835680751174
$ws.Range("F6").Value = 0.05633802816901409
$ws.Range("J6").Value = 0.2065727699530517
$ws.Range("O6").Value = 0.01408450704225352
$ws.Range("Q6").Value = 0.1455399061032864
$ws.Range("R6").Value = 0.07511737089201878
$ws.Range("S6").Value = 0.4131455399061033

# Row 7
$ws.Range("B7").Value = 0.06622516556291391
$ws.Range("D7").Value = 0.01986754966887417
$ws.Range("F7").Value = 0.03973509933774835
$ws.Range("J7").Value = 0.05298013245033113
$ws.Range("O7").Value = 0.01324503311258278
$ws.Range("Q7").Value = 0.1854304635761589
$ws.Range("R7").Value = 0.09271523178807947
$ws.Range("S7").Value = 0.5298013245033113

# Row 8
$ws.Range("B8").Value = 0.07423580786026202
$ws.Range("D8").Value = 0.03711790393013101
$ws.Range("E8").Value = 0.002183406113537118
$ws.Range("F8").Value = 0.07423580786026202
$ws.Range("J8").Value = 0.09606986899563319
$ws.Range("O8").Value = 0.008733624454148471
$ws.Range("Q8").Value = 0.2008733624454148
$ws.Range("R8").Value = 0.08733624454148471
$ws.Range("S8").Value = 0.4192139737991266

# Row 9
$ws.Range("B9").Value = 0.055
$ws.Range("D9").Value = 0.03
$ws.Range("E9").Value = 0.005
$ws.Range("F9").Value = 0.045
$ws.Range("J9").Value = 0.08
$ws.Range("O9").Value = 0.025
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.14
$ws.Range("S9").Value = 0.42

# Row 10
$ws.Range("B10").Value = 0.09614035087719298
$ws.Range("D10").Value = 0.01192982456140351
$ws.Range("E10").Value = 0.001403508771929824
$ws.Range("F10").Value = 0.06666666666666667
$ws.Range("J10").Value = 0.1080701754385965
$ws.Range("O10").Value = 0.01403508771929825
$ws.Range("Q10").Value = 0.2301754385964912
$ws.Range("R10").Value = 0.09964912280701754
$ws.Range("S10").Value = 0.3719298245614035

# Row 11
$ws.Range("G11").Value = 0.1074380165289256
$ws.Range("J11").Value = 0.1446280991735537
$ws.Range("K11").Value = 0.1983471074380165
$ws.Range("L11").Value = 0.5371900826446281
$ws.Range("S11").Value = 0.01239669421487603

# Row 12
$ws.Range("G12").Value = 0.7557251908396947
$ws.Range("J12").Value = 0.2137404580152672
$ws.Range("K12").Value = 0.007633587786259542
$ws.Range("L12").Value = 0.007633587786259542
$ws.Range("S12").Value = 0.01526717557251908

# Row 13
$ws.Range("G13").Value = 0.5306122448979592
$ws.Range("J13").Value = 0.4489795918367347
$ws.Range("S13").Value = 0.02040816326530612

# Row 15
$ws.Range("F15").Value = 0.02304147465437788
$ws.Range("H15").Value = 0.1797235023041475
$ws.Range("I15").Value = 0.05990783410138249
$ws.Range("J15").Value = 0.4101382488479263
$ws.Range("K15").Value = 0.03686635944700461
$ws.Range("M15").Value = 0.01382488479262673
$ws.Range("N15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.03686635944700461
$ws.Range("S15").Value = 0.2350230414746544

# Row 16
$ws.Range("F16").Value = 0.02395209580838323
$ws.Range("H16").Value = 0.1616766467065868
$ws.Range("I16").Value = 0.0718562874251497
$ws.Range("J16").Value = 0.4610778443113773
$ws.Range("K16").Value = 0.1197604790419162
$ws.Range("M16").Value = 0.005988023952095809
$ws.Range("O16").Value = 0.0658682634730539
$ws.Range("S16").Value = 0.08982035928143713

# Row 17
$ws.Range("F17").Value = 0.01160541586073501
$ws.Range("H17").Value = 0.1798839458413926
$ws.Range("I17").Value = 0.07543520309477757
$ws.Range("J17").Value = 0.4912959381044487
$ws.Range("K17").Value = 0.05996131528046422
$ws.Range("M17").Value = 0.01353965183752418
$ws.Range("N17").Value = 0.001934235976789168
$ws.Range("O17").Value = 0.05222437137330754
$ws.Range("S17").Value = 0.1141199226305609

# Row 18
$ws.Range("F18").Value = 0.01265822784810127
$ws.Range("H18").Value = 0.189873417721519
$ws.Range("I18").Value = 0.08438818565400844
$ws.Range("J18").Value = 0.4641350210970464
$ws.Range("K18").Value = 0.05907172995780591
$ws.Range("M18").Value = 0.01687763713080169
$ws.Range("N18").Value = 0.004219409282700422
$ws.Range("O18").Value = 0.07172995780590717
$ws.Range("S18").Value = 0.0970464135021097

# Row 19
$ws.Range("F19").Value = 0.01051051051051051
$ws.Range("H19").Value = 0.1921921921921922
$ws.Range("I19").Value = 0.0915915915915916
$ws.Range("J19").Value = 0.4091591591591592
$ws.Range("K19").Value = 0.08633633633633633
$ws.Range("M19").Value = 0.02627627627627628
$ws.Range("N19").Value = 0.0007507507507507507
$ws.Range("O19").Value = 0.06756756756756757
$ws.Range("S19").Value = 0.1156156156156156
